$d = $word.ActiveDocument

# Locate the paragraph that starts the block to be removed: "Ver no Jupiter ..."
# by searching the document text (robust to exact index positions).
$count = $d.Paragraphs.Count
$idxVerNoJupiter = -1
$idxCopyright = -1

for ($i = 1; $i -le $count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13)
    if ($txt -eq "Ver no Jupiter Salvar em pdf Salvar em docx") {
        $idxVerNoJupiter = $i
    }
    if ($txt -like "*Contact: luizeleno@usp.br*") {
        $idxCopyright = $i
    }
}

if ($idxVerNoJupiter -gt 0 -and $idxCopyright -gt 0) {
    # Also remove the blank paragraph immediately preceding "Ver no Jupiter ..."
    $startIndex = $idxVerNoJupiter
    $precedingPara = $d.Paragraphs.Item($idxVerNoJupiter - 1)
    if ($precedingPara.Range.Text.TrimEnd([char]13) -eq "") {
        $startIndex = $idxVerNoJupiter - 1
    }

    $startPara = $d.Paragraphs.Item($startIndex)
    $endPara = $d.Paragraphs.Item($idxCopyright)

    $rng = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $rng.Delete()
}
